$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.745.12"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.598.48"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'211.63"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "'0.0618"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'19.71"
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.823.27"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "1.628.76"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "'0.525"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "'64.99"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "26.702.74"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "0.0₃0742"
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("D19").Value = "'209.37"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("D24").Value = "'9.01"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").Value = "'144.58"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").Value = "'7.14"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'3.26"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").Value = "'2.98"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("E34").Value = "  +17.72%  "
$ws.Range("D35").Value = "1.274.70"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "'0.596"
$ws.Range("E38").Value = "  -3.78%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "'62.69"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Value = "1.734.43"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").Value = "'90.38"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").Value = "'7.50"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("E51").Value = "  +0.14%  "
